# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2  = 3
    3  = 3
    4  = 3
    5  = 0
    6  = 3
    7  = 3
    8  = 1
    9  = 6
    10 = 0
    11 = 3
    12 = 1
    14 = 2
    17 = 0
    18 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
